# Commit: "Publication des derniers chapitres"
#
# The Pokemon names in column B that carried French accented characters
# are rewritten without diacritics (e.g. "Salamèche" -> "Salameche").
# This is exactly the set of `pokemon` sheet "nom" cells whose accented
# text is being replaced by its accent-stripped form.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pokemon")

$renames = [ordered]@{
    "B5"   = "Salameche"
    "B36"  = "Melofee"
    "B37"  = "Melodelfe"
    "B50"  = "Aeromite"
    "B57"  = "Ferosinge"
    "B62"  = "Tetarte"
    "B70"  = "Chetiflor"
    "B82"  = "Magneti"
    "B83"  = "Magneton"
    "B102" = "electrode"
    "B113" = "Rhinoferos"
    "B118" = "Hypocean"
    "B119" = "Poissirene"
    "B124" = "Insecateur"
    "B126" = "elektek"
    "B131" = "Leviator"
    "B133" = "Metamorph"
    "B134" = "evoli"
    "B143" = "Ptera"
    "B146" = "electhor"
    "B155" = "Meganium"
    "B156" = "Hericendre"
    "B174" = "Melo"
    "B193" = "Heliatronc"
    "B199" = "Cornebre"
    "B201" = "Feuforeve"
    "B203" = "Qulbutoke"
    "B224" = "Remoraid"
    "B227" = "Demanta"
    "B230" = "Demolosse"
    "B240" = "elekid"
    "B242" = "ecremeuh"
}

foreach ($cellRef in $renames.Keys) {
    $ws.Range($cellRef).Value = $renames[$cellRef]
}

# Move the saved view position to the bottom of the list, matching the
# author scrolling down to review the newly-edited rows before publishing.
$ws.Activate()
$ws.Range("B256").Select()

Write-Output "Renamed $($renames.Count) pokemon names to strip accents."
